# Append the next daily mod-count reading for 逃离鸭科夫 as row 86
# (Date=2026/02/04, Game=逃离鸭科夫, ModCount=1171), matching the
# formatting of the preceding data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data rows use center/center alignment with the default General
# number format (style index 1 in this sheet). Apply the same alignment to
# the new row before filling in values.
$ws.Range("A86:C86").HorizontalAlignment = -4108
$ws.Range("A86:C86").VerticalAlignment = -4108

# Assigning the date-like text straight to Range.Value would let Excel's
# smart entry reinterpret "2026/02/04" as a real date serial, changing both
# the cell's type and its number format. To keep it as plain text (as the
# source data already is), stage the literal string as a text-producing
# formula in a scratch cell, copy it, and paste only the *value* into the
# target cell -- this preserves the destination's existing type/style while
# still landing the literal text "2026/02/04".
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""2026/02/04"""
$scratch.Copy()
$ws.Range("A86").PasteSpecial(-4163)
$scratch.Value = ""

$ws.Range("B86").Value = "逃离鸭科夫"
$ws.Range("C86").Value = 1171
